# Updates FFXIV market-data columns (H..N) across the 8 job sheets,
# reflecting the latest scheduled market-board price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2322.4075
$ws.Range("I19").Value = 3958.0715
$ws.Range("K19").Value = 3958.0715
$ws.Range("M19").Value = -3783.0715

$ws.Range("H33").Value = 803.2973
$ws.Range("I33").Value = 582.1739
$ws.Range("J33").Value = 1166.5714
$ws.Range("K33").Value = 582.1739
$ws.Range("L33").Value = 1166.5714
$ws.Range("M33").Value = -353.1739
$ws.Range("N33").Value = -1624.5714

$ws.Range("H98").Value = 1481.5
$ws.Range("I98").Value = 1142
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 1142
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 356
$ws.Range("N98").Value = -5496

$ws.Range("H113").Value = 3484.9375
$ws.Range("I113").Value = 2782.375
$ws.Range("J113").Value = 4187.5
$ws.Range("K113").Value = 2782.375
$ws.Range("L113").Value = 4187.5
$ws.Range("M113").Value = 471.625
$ws.Range("N113").Value = -10695.5

$ws.Range("H122").Value = 1481.5
$ws.Range("I122").Value = 1142
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3426
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -976
$ws.Range("N122").Value = -12400

$ws.Range("H138").Value = 2133.6448
$ws.Range("I138").Value = 1374.4762
$ws.Range("J138").Value = 2423.509
$ws.Range("K138").Value = 4123.4286
$ws.Range("L138").Value = 7270.527
$ws.Range("M138").Value = 1016.5714
$ws.Range("N138").Value = -17550.527

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700
$ws.Range("I2").Value = 500
$ws.Range("K2").Value = 500
$ws.Range("M2").Value = -387

$ws.Range("H34").Value = 15000
$ws.Range("J34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15542

$ws.Range("H42").Value = 6031
$ws.Range("J42").Value = 6031
$ws.Range("L42").Value = 6031
$ws.Range("N42").Value = -7003

$ws.Range("H61").Value = 3032.745
$ws.Range("I61").Value = 1962.2667
$ws.Range("J61").Value = 4562
$ws.Range("K61").Value = 1962.2667
$ws.Range("L61").Value = 4562
$ws.Range("M61").Value = -1750.2667
$ws.Range("N61").Value = -4986

$ws.Range("H63").Value = 4024.6365
$ws.Range("I63").Value = 4168.6787
$ws.Range("J63").Value = 3218
$ws.Range("K63").Value = 4168.6787
$ws.Range("L63").Value = 3218
$ws.Range("M63").Value = -3482.6787
$ws.Range("N63").Value = -4590

$ws.Range("H66").Value = 4024.6365
$ws.Range("I66").Value = 4168.6787
$ws.Range("J66").Value = 3218
$ws.Range("K66").Value = 20843.3935
$ws.Range("L66").Value = 16090
$ws.Range("M66").Value = -17411.3935
$ws.Range("N66").Value = -22954

$ws.Range("H116").Value = 700
$ws.Range("I116").Value = 500
$ws.Range("K116").Value = 500
$ws.Range("M116").Value = 1794

$ws.Range("H132").Value = 3041.7817
$ws.Range("I132").Value = 2979.0557
$ws.Range("K132").Value = 8937.167099999999
$ws.Range("M132").Value = -6407.167099999999

$ws.Range("H136").Value = 3032.745
$ws.Range("I136").Value = 1962.2667
$ws.Range("J136").Value = 4562
$ws.Range("K136").Value = 5886.800099999999
$ws.Range("L136").Value = 13686
$ws.Range("M136").Value = -3336.800099999999
$ws.Range("N136").Value = -18786

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -386

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()  # profit no longer computable; cell removed

$ws.Range("H80").Value = 343.96
$ws.Range("J80").Value = 406.5625
$ws.Range("L80").Value = 406.5625
$ws.Range("N80").Value = -2402.5625

$ws.Range("H83").Value = 343.96
$ws.Range("J83").Value = 406.5625
$ws.Range("L83").Value = 2032.8125
$ws.Range("N83").Value = -12016.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2780
$ws.Range("I16").Value = 3250
$ws.Range("K16").Value = 3250
$ws.Range("M16").Value = -2963

$ws.Range("H31").Value = 3179.4314
$ws.Range("I31").Value = 2213.182
$ws.Range("J31").Value = 4950.8887
$ws.Range("K31").Value = 2213.182
$ws.Range("L31").Value = 4950.8887
$ws.Range("M31").Value = -1918.182
$ws.Range("N31").Value = -5540.8887

$ws.Range("H34").Value = 3179.4314
$ws.Range("I34").Value = 2213.182
$ws.Range("J34").Value = 4950.8887
$ws.Range("K34").Value = 2213.182
$ws.Range("L34").Value = 4950.8887
$ws.Range("M34").Value = -2011.182
$ws.Range("N34").Value = -5354.8887

$ws.Range("H35").Value = 83335490
$ws.Range("I35").Value = 125000730
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 125000730
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -125000436
$ws.Range("N35").Value = -5588

$ws.Range("H113").Value = 2780
$ws.Range("I113").Value = 3250
$ws.Range("K113").Value = 3250
$ws.Range("M113").Value = -1080

$ws.Range("H132").Value = 2361.3333
$ws.Range("I132").Value = 1321.2778
$ws.Range("J132").Value = 3401.389
$ws.Range("K132").Value = 3963.8334
$ws.Range("L132").Value = 10204.167
$ws.Range("M132").Value = -1433.8334
$ws.Range("N132").Value = -15264.167

$ws.Range("H134").Value = 1975.9412
$ws.Range("I134").Value = 1158.2941
$ws.Range("K134").Value = 3474.8823
$ws.Range("M134").Value = -939.8823000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38.52381
$ws.Range("I12").Value = 27.125
$ws.Range("K12").Value = 81.375
$ws.Range("M12").Value = 91.625

$ws.Range("H23").Value = 504.68
$ws.Range("J23").Value = 524
$ws.Range("L23").Value = 1572
$ws.Range("N23").Value = -2042

$ws.Range("H92").Value = 949.5
$ws.Range("J92").Value = 966
$ws.Range("L92").Value = 2898
$ws.Range("N92").Value = -5394

$ws.Range("H116").Value = 2806.889
$ws.Range("I116").Value = 646
$ws.Range("K116").Value = 1938
$ws.Range("M116").Value = 1504

$ws.Range("H120").Value = 10460.363
$ws.Range("I120").Value = 9600
$ws.Range("K120").Value = 28800
$ws.Range("M120").Value = -23962

$ws.Range("H131").Value = 2309.238
$ws.Range("J131").Value = 1861.4706
$ws.Range("L131").Value = 5584.4118
$ws.Range("N131").Value = -15664.4118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724

$ws.Range("H102").Value = 3238.7856
$ws.Range("I102").Value = 3438.88
$ws.Range("J102").Value = 1571.3334
$ws.Range("K102").Value = 3438.88
$ws.Range("L102").Value = 1571.3334
$ws.Range("M102").Value = -1816.88
$ws.Range("N102").Value = -4815.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3500.611
$ws.Range("I40").Value = 3262.4614
$ws.Range("J40").Value = 4119.8
$ws.Range("K40").Value = 3262.4614
$ws.Range("L40").Value = 4119.8
$ws.Range("M40").Value = -3126.4614
$ws.Range("N40").Value = -4391.8

$ws.Range("H48").Value = 3379.3333
$ws.Range("I48").Value = 919
$ws.Range("J48").Value = 8300
$ws.Range("K48").Value = 919
$ws.Range("L48").Value = 8300
$ws.Range("M48").Value = -258
$ws.Range("N48").Value = -9622

$ws.Range("H82").Value = 2357.3125
$ws.Range("I82").Value = 2282.913
$ws.Range("J82").Value = 2547.4443
$ws.Range("K82").Value = 2282.913
$ws.Range("L82").Value = 2547.4443
$ws.Range("M82").Value = -1921.913
$ws.Range("N82").Value = -3269.4443

$ws.Range("H85").Value = 2357.3125
$ws.Range("I85").Value = 2282.913
$ws.Range("J85").Value = 2547.4443
$ws.Range("K85").Value = 2282.913
$ws.Range("L85").Value = 2547.4443
$ws.Range("M85").Value = -1034.913
$ws.Range("N85").Value = -5043.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10826

$ws.Range("H42").Value = 9000
$ws.Range("J42").Value = 9000
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -9756

$ws.Range("H43").Value = 12860.333
$ws.Range("I43").Value = 10001
$ws.Range("J43").Value = 14290
$ws.Range("K43").Value = 10001
$ws.Range("L43").Value = 14290
$ws.Range("M43").Value = -9852
$ws.Range("N43").Value = -14588

$ws.Range("H113").Value = 323.75
$ws.Range("I113").Value = 309.05264
$ws.Range("J113").Value = 603
$ws.Range("K113").Value = 927.15792
$ws.Range("L113").Value = 1809
$ws.Range("M113").Value = 1242.84208
$ws.Range("N113").Value = -6149

$ws.Range("H136").Value = 14494318
$ws.Range("I136").Value = 25001186
$ws.Range("K136").Value = 75003558
$ws.Range("M136").Value = -75001008
